# Update "paises" (countries) workbook - refresh COVID-19 data snapshot
# (new case counts for several countries) and re-rank the countries whose
# totals changed order relative to their neighbours in the list.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update "last updated" timestamp banner in A1 ---
$ws.Range("A1").Value = "Datos actualizados a 21 de Abril de 2020 a las 11:52"

# --- Row 14: Belgica - refresh active cases / recovered / critical cases ---
$ws.Range("D14").Value = 9002
$ws.Range("E14").Value = 25956
$ws.Range("F14").Value = 1079

# --- Row 47: Malasia - refresh full stat line ---
$ws.Range("B47").Value = 5482
$ws.Range("C47").Value = 57
$ws.Range("D47").Value = 3349
$ws.Range("F47").Value = 43
$ws.Range("G47").Value = 3
$ws.Range("H47").Value = 92

# --- Rows 50-51: Finlandia overtakes Colombia in ranking ---
# Row 50 becomes Finlandia (updated data), row 51 becomes Colombia (previous Colombia data)
$ws.Range("A50").Value = "Finlandia"
$ws.Range("B50").Value = 4014
$ws.Range("C50").Value = 146
$ws.Range("D50").Value = 2000
$ws.Range("E50").Value = 1916
$ws.Range("F50").Value = 67
$ws.Range("G50").Value = 0
$ws.Range("H50").Value = 98

$ws.Range("A51").Value = "Colombia"
$ws.Range("B51").Value = 3977
$ws.Range("C51").Value = 0
$ws.Range("D51").Value = 804
$ws.Range("E51").Value = 2984
$ws.Range("F51").Value = 98
$ws.Range("G51").Value = 0
$ws.Range("H51").Value = 189

# --- Row 63: Moldavia - refresh full stat line ---
$ws.Range("B63").Value = 2080
$ws.Range("C63").Value = 85
$ws.Range("D63").Value = 412
$ws.Range("E63").Value = 1657
$ws.Range("F63").Value = 46
$ws.Range("G63").Value = 2
$ws.Range("H63").Value = 11

# --- Rows 98-100: Albania overtakes Bolivia and Kirguistan in ranking ---
# Row 98 becomes Albania (updated data), row 99 becomes Bolivia (previous data),
# row 100 becomes Kirguistan (previous data)
$ws.Range("A98").Value = "Albania"
$ws.Range("B98").Value = 609
$ws.Range("C98").Value = 25
$ws.Range("D98").Value = 345
$ws.Range("E98").Value = 238
$ws.Range("F98").Value = 10
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 26

$ws.Range("A99").Value = "Bolivia"
$ws.Range("B99").Value = 598
$ws.Range("C99").Value = 34
$ws.Range("D99").Value = 37
$ws.Range("E99").Value = 527
$ws.Range("F99").Value = 3
$ws.Range("G99").Value = 1
$ws.Range("H99").Value = 34

$ws.Range("A100").Value = "Kirguistan"
$ws.Range("B100").Value = 590
$ws.Range("C100").Value = 22
$ws.Range("D100").Value = 216
$ws.Range("E100").Value = 367
$ws.Range("F100").Value = 5
$ws.Range("G100").Value = 0
$ws.Range("H100").Value = 7
